$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UBER")

# Row 12: Accounts Payable
$ws.Range("B12").Value = 235000000.0
$ws.Range("C12").Value = 240000000.0
$ws.Range("D12").Value = 253000000.0
$ws.Range("E12").Value = 215000000.0
$ws.Range("F12").Value = 272000000.0

# Row 21: Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = 818000000.0
$ws.Range("C21").Value = 787000000.0
$ws.Range("D21").Value = 752000000.0
$ws.Range("E21").Value = 793000000.0
$ws.Range("F21").Value = 1027000000.0
